# perbaikan view pada master barang sesuai dengan struktur terbaru
# Adjust the "master barang" export sheet: row 4 (IN=24/OUT="-") becomes the
# template for three extra "-" rows, expanding the table from A1:E4 to A1:E7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure gridlines / row-col headers stay visible (matches original view).
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.DisplayHeadings = $true

# Copy the formatting (borders/alignment/style) of the existing data row (row 4)
# down into the three new rows (5-7) before writing any values into them.
$ws.Range("A4:E4").Copy()
$ws.Range("A5:E7").PasteSpecial(-4122)  # xlPasteFormats

# Row 4: keep the product code / stock, but now IN = 24, OUT = "-", Harga = 12000
$ws.Range("C4").Value = 24
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = 12000

# Row 5: re-insert what used to be row 4 (IN = "-", OUT = 1, Harga = 7000)
$ws.Range("A5").Value = "65c5c56b6807b"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = "-"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 7000

# Row 6: another IN = 24 / OUT = "-" row, Harga = 12000
$ws.Range("A6").Value = "65c5c56b6807b"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 24
$ws.Range("D6").Value = "-"
$ws.Range("E6").Value = 12000

# Row 7: another IN = 24 / OUT = "-" row, Harga = 200000
$ws.Range("A7").Value = "65c5c56b6807b"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 24
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = 200000

# Match the updated selection/active cell from the diff (dimension grows to E7).
$ws.Range("E7").Select() | Out-Null
